$wb = $excel.ActiveWorkbook

# Insert two new sheets after "scenario" and before "interventions":
#   region_definitions
#   scenario_definitions
$scenarioSheet = $wb.Worksheets.Item("scenario")

$regionDefs = $wb.Worksheets.Add($null, $scenarioSheet)
$regionDefs.Name = "region_definitions"
$regionDefs.Range("A1").Value = "name"
$regionDefs.Range("B1").Value = "description"
$regionDefs.Range("C1").Value = "filename"

$scenarioDefs = $wb.Worksheets.Add($null, $regionDefs)
$scenarioDefs.Name = "scenario_definitions"
$scenarioDefs.Range("A1").Value = "name"
$scenarioDefs.Range("B1").Value = "description"
$scenarioDefs.Range("C1").Value = "filename"

# Set each sheet's selection to match the saved view state.
$regionDefs.Activate()
$regionDefs.Range("D1").Select()

$scenarioDefs.Activate()
$scenarioDefs.Range("B2").Select()

# The "interventions" sheet becomes the active tab.
$interventions = $wb.Worksheets.Item("interventions")
$interventions.Activate()
$interventions.Range("E22").Select()
